$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "value" column header to "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# The series grew from 52 to 83 quarterly observations (rows 2-84).
# Extend the date-formatted style from A53 down through the new rows so the
# appended dates keep the same YYYY-MM-DD HH:MM:SS number format.
$ws.Range("A53").Copy($ws.Range("A54:A84"))

# Refresh the full data series (dates are Excel serial numbers, values are the
# quarter-on-quarter percentage changes for the first-release vintage).
$ws.Range("A2").Value = 38398
$ws.Range("B2").Value = 1.3
$ws.Range("A3").Value = 38487
$ws.Range("B3").Value = -0.5
$ws.Range("A4").Value = 38579
$ws.Range("B4").Value = 0.09999999999999987
$ws.Range("A5").Value = 38671
$ws.Range("B5").Value = -0.4
$ws.Range("A6").Value = 38763
$ws.Range("B6").Value = 0.7
$ws.Range("A7").Value = 38852
$ws.Range("B7").Value = 0.5
$ws.Range("A8").Value = 38944
$ws.Range("B8").Value = 0
$ws.Range("A9").Value = 39036
$ws.Range("B9").Value = 1.9
$ws.Range("A10").Value = 39128
$ws.Range("B10").Value = -1
$ws.Range("A11").Value = 39217
$ws.Range("B11").Value = 1.1
$ws.Range("A12").Value = 39309
$ws.Range("B12").Value = -0.2
$ws.Range("A13").Value = 39401
$ws.Range("B13").Value = 0.5
$ws.Range("A14").Value = 39493
$ws.Range("B14").Value = -0.2
$ws.Range("A15").Value = 39583
$ws.Range("B15").Value = 0.5
$ws.Range("A16").Value = 39675
$ws.Range("B16").Value = -1.7
$ws.Range("A17").Value = 39767
$ws.Range("B17").Value = -2
$ws.Range("A18").Value = 39859
$ws.Range("B18").Value = -2.168616218514021
$ws.Range("A19").Value = 39948
$ws.Range("B19").Value = 1.9
$ws.Range("A20").Value = 40040
$ws.Range("B20").Value = -0.5
$ws.Range("A21").Value = 40132
$ws.Range("B21").Value = 1.6
$ws.Range("A22").Value = 40224
$ws.Range("B22").Value = -1.1
$ws.Range("A23").Value = 40313
$ws.Range("B23").Value = 0.8
$ws.Range("A24").Value = 40405
$ws.Range("B24").Value = 0.2696221576284057
$ws.Range("A25").Value = 40497
$ws.Range("B25").Value = 0.5
$ws.Range("A26").Value = 40589
$ws.Range("B26").Value = 0.5
$ws.Range("A27").Value = 40678
$ws.Range("B27").Value = 0.1000000000000001
$ws.Range("A28").Value = 40770
$ws.Range("B28").Value = 0.09999999999999987
$ws.Range("A29").Value = 40862
$ws.Range("B29").Value = -0.4
$ws.Range("A30").Value = 40954
$ws.Range("B30").Value = 0.9
$ws.Range("A31").Value = 41044
$ws.Range("B31").Value = 0.7
$ws.Range("A32").Value = 41136
$ws.Range("B32").Value = 0.3
$ws.Range("A33").Value = 41228
$ws.Range("B33").Value = -0.7000000000000001
$ws.Range("A34").Value = 41320
$ws.Range("B34").Value = 0.09999999999999998
$ws.Range("A35").Value = 41409
$ws.Range("B35").Value = 0.2999999999999999
$ws.Range("A36").Value = 41501
$ws.Range("B36").Value = -0.4
$ws.Range("A37").Value = 41593
$ws.Range("B37").Value = 0.7000000000000001
$ws.Range("A38").Value = 41685
$ws.Range("B38").Value = -0.9
$ws.Range("A39").Value = 41774
$ws.Range("B39").Value = 0
$ws.Range("A40").Value = 41866
$ws.Range("B40").Value = 0.2000000000000001
$ws.Range("A41").Value = 41958
$ws.Range("B41").Value = -0.3
$ws.Range("A42").Value = 42050
$ws.Range("B42").Value = -0.2
$ws.Range("A43").Value = 42139
$ws.Range("B43").Value = 0.6
$ws.Range("A44").Value = 42231
$ws.Range("B44").Value = -0.4
$ws.Range("A45").Value = 42323
$ws.Range("B45").Value = -0.5
$ws.Range("A46").Value = 42415
$ws.Range("B46").Value = -0.1
$ws.Range("A47").Value = 42505
$ws.Range("B47").Value = 0.5
$ws.Range("A48").Value = 42597
$ws.Range("B48").Value = -0.3
$ws.Range("A49").Value = 42689
$ws.Range("B49").Value = -0.2
$ws.Range("A50").Value = 42781
$ws.Range("B50").Value = 0.4
$ws.Range("A51").Value = 42870
$ws.Range("B51").Value = -0.4
$ws.Range("A52").Value = 42962
$ws.Range("B52").Value = 0.4
$ws.Range("A53").Value = 43054
$ws.Range("B53").Value = 0.5
$ws.Range("A54").Value = 43146
$ws.Range("B54").Value = -0.1
$ws.Range("A55").Value = 43235
$ws.Range("B55").Value = -0.4
$ws.Range("A56").Value = 43327
$ws.Range("B56").Value = -1
$ws.Range("A57").Value = 43419
$ws.Range("B57").Value = 0
$ws.Range("A58").Value = 43511
$ws.Range("B58").Value = 0.2
$ws.Range("A59").Value = 43600
$ws.Range("B59").Value = -0.5
$ws.Range("A60").Value = 43692
$ws.Range("B60").Value = 0.5
$ws.Range("A61").Value = 43784
$ws.Range("B61").Value = -0.6
$ws.Range("A62").Value = 43876
$ws.Range("B62").Value = -0.8
$ws.Range("A63").Value = 43966
$ws.Range("B63").Value = -2.8
$ws.Range("A64").Value = 44058
$ws.Range("B64").Value = 2.03040224812923
$ws.Range("A65").Value = 44150
$ws.Range("B65").Value = 0.4339050167294337
$ws.Range("A66").Value = 44242
$ws.Range("B66").Value = 1.197090680270919
$ws.Range("A67").Value = 44331
$ws.Range("B67").Value = -0.6441697788099752
$ws.Range("A68").Value = 44423
$ws.Range("B68").Value = 0.4901411037136303
$ws.Range("A69").Value = 44515
$ws.Range("B69").Value = -0.09098140646410988
$ws.Range("A70").Value = 44607
$ws.Range("B70").Value = -0.3952916234765647
$ws.Range("A71").Value = 44696
$ws.Range("B71").Value = -0.2380782088493735
$ws.Range("A72").Value = 44788
$ws.Range("B72").Value = -0.02511582766690132
$ws.Range("A73").Value = 44880
$ws.Range("B73").Value = 0.6876823391013496
$ws.Range("A74").Value = 44972
$ws.Range("B74").Value = 0.3720617293507145
$ws.Range("A75").Value = 45061
$ws.Range("B75").Value = -0.1123255314657629
$ws.Range("A76").Value = 45153
$ws.Range("B76").Value = -0.1387895598915543
$ws.Range("A77").Value = 45245
$ws.Range("B77").Value = 0.1213803088128225
$ws.Range("A78").Value = 45337
$ws.Range("B78").Value = -0.07951810869463416
$ws.Range("A79").Value = 45427
$ws.Range("B79").Value = 0.3863682696630121
$ws.Range("A80").Value = 45519
$ws.Range("B80").Value = -0.006123215295980228
$ws.Range("A81").Value = 45611
$ws.Range("B81").Value = 0.04342916022020096
$ws.Range("A82").Value = 45703
$ws.Range("B82").Value = -0.2736421272901388
$ws.Range("A83").Value = 45792
$ws.Range("B83").Value = -0.1489623566660376
$ws.Range("A84").Value = 45884
$ws.Range("B84").Value = -0.1051295019602515
